$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the new
# header cells so they match the other bold/bordered header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record (Wins/Losses/Ties) is the same for every player row.
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = 93
    $ws.Cells.Item($row, 31).Value = 69
    $ws.Cells.Item($row, 32).Value = 0
}
